# Update the "Plants" variables table:
#   - insert a new "Initial Value" column between "Tune" and "Estimated sensetivity"
#   - fill in the initial-value numbers for the two existing variables
#   - add four more Cryptomonad variable rows (EMort, MinLightSat, TempRespSlope...)
#   - carry the row-height / column-width cosmetics that Excel's autofit produced

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column F ("Initial Value") before the old "Estimated sensetivity"
#     column, which slides from F to G and keeps its own width/content. ---
$ws.Columns.Item(6).Insert()

$ws.Range("F1").Value = "Initial Value"
$ws.Range("F2").Value = 3.9
$ws.Range("F3").Value = 0.07

# --- Column widths (character-width units; Excel stores width+5/6 in the xml) ---
$ws.Columns.Item(1).ColumnWidth = 15.166666666666666
$ws.Columns.Item(2).ColumnWidth = 25.022135416666668
$ws.Columns.Item(3).ColumnWidth = 25.022135416666668
$ws.Columns.Item(4).ColumnWidth = 17.022135416666668
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668

# --- Insert three more rows right after row 3 for the newly-tracked variables ---
$ws.Rows.Item(4).Resize(3).Insert()

# Row 4: Exponential factor for suboptimal conditions -> "EMort"
$ws.Range("A4").Value = "Cryptomonad"
$ws.Range("B4").Value = "OtherAlg2: [Cryptomonad]"
$ws.Range("C4").Value = "Exponential factor for suboptimal conditions"
$ws.Range("D4").Value = """EMort"""
$ws.Range("E4").Value = "yes"
$ws.Range("F4").Value = 0.09

# Row 5: Minimum light saturation -> "MinLightSat"
$ws.Range("A5").Value = "Cryptomonad"
$ws.Range("B5").Value = "OtherAlg2: [Cryptomonad]"
$ws.Range("C5").Value = "Minimum light saturation"
$ws.Range("D5").Value = """MinLightSat"""
$ws.Range("E5").Value = "yes"
$ws.Range("F5").Value = 11

# Row 6: Slope or rate of change per 10 degC temperature change -> "TempRespSlope"
$ws.Range("A6").Value = "Cryptomonad"
$ws.Range("B6").Value = "OtherAlg2: [Cryptomonad]"
$ws.Range("C6").Value = "Slope or rate of change per 10°C temperature change (maybe photosynthesis?)"
$ws.Range("D6").Value = """TempRespSlope"""
$ws.Range("E6").Value = "yes"
$ws.Range("F6").Value = 2

# --- Row heights matching the wrapped-text autofit result ---
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 45

# --- Selection left where the editor ended up ---
$ws.Range("G10").Select() | Out-Null
